$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching the style of the existing header row (F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New data cells for row 2
$ws.Range("G2").Value = 0.125854933266722
$ws.Range("H2").Value = 0.9890000000000001
